$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price column D, volume/change column E)
# Number-like price values are prefixed with a leading apostrophe so Excel
# stores them as literal text (matching the source inlineStr cells) instead
# of auto-converting them to numeric values and losing formatting such as
# trailing zeros (e.g. "0.05110") or thousands-grouped strings.

$ws.Range("D2").Value = "26.086.47"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.647.77"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'218.13"
$ws.Range("D6").Value = "'0.5199"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "'0.2624"
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("D9").Value = "'0.06302"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("D10").Value = "'20.31"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("D11").Value = "'0.07671"
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("D12").Value = "'4.586"
$ws.Range("E12").Value = "  +2.43%  "
$ws.Range("D13").Value = "1.632.09"
$ws.Range("E13").Value = "  -3.80%  "
$ws.Range("D14").Value = "1.875.57"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").Value = "0.0₅8109"
$ws.Range("E16").Value = "  +1.47%  "
$ws.Range("D17").Value = "'65.09"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").Value = "26.066.06"
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Value = "'4.596"
$ws.Range("E20").Value = "  -1.00%  "
$ws.Range("D21").Value = "'192.85"
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("D22").Value = "'10.42"
$ws.Range("E22").Value = "  +3.18%  "
$ws.Range("D23").Value = "'5.908"
$ws.Range("E23").Value = "  -0.70%  "
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").Value = "'144.15"
$ws.Range("E25").Value = "  -1.46%  "
$ws.Range("E26").Value = "  -1.69%  "
$ws.Range("D27").Value = "'7.178"
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("D28").Value = "'15.84"
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("E29").Value = "  +1.77%  "
$ws.Range("D30").Value = "'0.05355"
$ws.Range("E30").Value = "  -4.65%  "
$ws.Range("D31").Value = "'1.269"
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("D32").Value = "'3.442"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").Value = "'3.322"
$ws.Range("E33").Value = "  -1.35%  "
$ws.Range("D34").Value = "'1.548"
$ws.Range("E34").Value = "  -3.04%  "
$ws.Range("D35").Value = "'2.416"
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("D36").Value = "'2.780"
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("D37").Value = "'0.9414"
$ws.Range("D38").Value = "'0.5592"
$ws.Range("E38").Value = "  -1.13%  "
$ws.Range("D39").Value = "'0.01571"
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("D40").Value = "'5.778"
$ws.Range("E40").Value = "  -2.98%  "
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").Value = "1.026.60"
$ws.Range("E42").Value = "  -2.48%  "
$ws.Range("D43").Value = "'0.8257"
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("E44").Value = "  -1.29%  "
$ws.Range("D45").Value = "1.785.84"
$ws.Range("E46").Value = "  +10.80%  "
$ws.Range("D47").Value = "'57.29"
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("E48").Value = "  -0.52%  "
$ws.Range("D49").Value = "'0.4316"
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("D50").Value = "'7.904"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("D51").Value = "'0.05110"
$ws.Range("E51").Value = "  -3.94%  "
